$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.176.50"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").Value = "1.901.55"
$ws.Range("E3").Value = "  +1.96%  "

$ws.Range("E4").Value = "  +1.05%  "

$ws.Range("D5").Value = "'337.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").Value = "'1.015"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("D7").Value = "'0.4829"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.92%  "

$ws.Range("D8").Value = "'0.3998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.66%  "

$ws.Range("D9").Value = "'46.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.28%  "

$ws.Range("D10").Value = "'0.08106"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.49%  "

$ws.Range("D11").Value = "'1.031"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "'22.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.08%  "

$ws.Range("D13").Value = "1.902.56"
$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("D14").Value = "'6.059"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("D15").Value = "'7.277"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("D16").Value = "'1.017"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'88.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.74%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06786"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.89%  "

$ws.Range("D19").Value = "'0.00001055"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("D20").Value = "'17.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").Value = "'1.014"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("D22").Value = "28.163.20"
$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("D23").Value = "'5.572"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.76%  "

$ws.Range("D24").Value = "'11.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").Value = "'2.361"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("D26").Value = "2.126.67"
$ws.Range("E26").Value = "  +2.05%  "

$ws.Range("D27").Value = "'161.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").Value = "'20.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "'2.139"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("D30").Value = "'5.619"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").Value = "'122.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.9889"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.64%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.09676"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.96%  "

$ws.Range("D34").Value = "'3.651"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("D35").Value = "'5.403"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.93%  "

$ws.Range("D36").Value = "'1.381"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.24%  "

$ws.Range("D37").Value = "'0.06148"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.61%  "

$ws.Range("D38").Value = "'0.02272"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").Value = "'1.217"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("D40").Value = "'8.312"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "'1.014"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.92%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6043"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.24%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1915"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.40%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'10.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.55%  "

$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5733"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").Value = "'1.960"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.45%  "

$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "'3.392"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06846"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'113.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
